# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "36.653.60"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +2.31%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "2.104.91"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +11.16%  "
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "249.41"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.99%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.671"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -2.93%  "
$ws.Cells.Item(7, 5).Value = "  -0.03%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "45.49"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +4.81%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "61.21"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +7.36%  "
$ws.Cells.Item(10, 5).Value = "  +1.92%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.0731"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -3.50%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.0992"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.55%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "14.65"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -1.73%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "2.402.42"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +10.74%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.845"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +6.72%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "2.094.57"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +10.53%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "5.05"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.47%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "36.662.37"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +2.27%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "72.88"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.58%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0821"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -1.21%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "241.04"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -2.72%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "12.89"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.33%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "5.07"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -2.35%  "
$ws.Cells.Item(24, 5).Value = "  +0.19%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "2.47"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -9.88%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "170.06"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.78%  "
$ws.Cells.Item(27, 2).Value = "Cosmos"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "9.01"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +3.52%  "
$ws.Cells.Item(28, 2).Value = "EthereumClassic"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "20.48"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +10.74%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "2.00"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -8.23%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "0.123"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -4.39%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "22.45"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +50.76%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "4.44"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.32%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.0597"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -2.07%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "0.0905"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +15.56%  "
$ws.Cells.Item(35, 5).Value = "  -0.14%  "
$ws.Cells.Item(36, 5).Value = "  -0.84%  "
$ws.Cells.Item(37, 5).Value = "  +18.57%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "4.08"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -4.92%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.908"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +5.62%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "1.36"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -8.57%  "
$ws.Cells.Item(41, 5).Value = "  +8.84%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "99.82"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.06%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.0219"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -4.14%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "2.81"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +16.84%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "16.29"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -3.95%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "1.362.44"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +3.57%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.0835"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +3.08%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "2.290.14"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +10.59%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "2.84"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +3.03%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "2.27"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -3.04%  "
$ws.Cells.Item(51, 5).Value = "  +16.63%  "
